$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "score" column header in I1
$ws.Range("I1").Value = "score"

# Row 3 (KIMIT YOUSSEF): trim the Experience professionelle (F3) down to the
# internship line only, dropping the project description sentence.
$ws.Range("F3").Value = "Stagiaire Développeur chez GSDEV, Derb-Omar, Casablanca,`nMaroc.,"

# Row 3: trim Formation (G3) down to just the diploma/cycle titles before the
# colon/period, dropping the detailed descriptions.
$ws.Range("G3").Value = "Cycle d’ingénieur :,Cycle d’ingénieur :,Diplôme d'études universitaires générales :,Diplôme de Baccalauréat en Science physiques,"

# Row 3: trim Competence (H3) down to the core tech skill list ending in
# "Unity,", dropping the long free-text project/certification descriptions.
$ws.Range("H3").Value = "Java,Python,R,C,C#,PHP,Javascript,Angular,React,Spring Boot,Hibernate,Bootstrap,NumPy,Pandas,Matplotlib,Seaborn,Scikit-Learn,MySQL,Oracle,MongoDB,Unity,"

# The embedded line break in F3 makes Excel auto-expand the row height;
# restore the row to its natural auto-fit height so no explicit height is
# persisted (matching the original layout).
$ws.Rows.Item(3).EntireRow.AutoFit()
